$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates from the diff: Price (D) and Volume(1h) (E) columns.
# NumberFormat is forced to text ("@") before the assignment so that
# numeric-looking / percent-looking strings are preserved as literal
# text (matching the source workbook's inlineStr cells) instead of
# being auto-converted to numbers by Excel's input parser.
$updates = @{
    "D2" = "307.34"
    "E2" = "1.30%"
    "D3" = "35.88"
    "E3" = "2.27%"
    "D4" = "5.110"
    "E4" = "0.88%"
    "D5" = "0.08086"
    "E5" = "1.06%"
    "D6" = "1.952"
    "E6" = "1.13%"
    "D7" = "4.207"
    "E7" = "3.71%"
    "E8" = "0.14%"
    "D9" = "0.9295"
    "E9" = "0.97%"
    "D10" = "0.1383"
    "E10" = "12.50%"
    "D11" = "0.1925"
    "E11" = "4.24%"
    "D12" = "0.09227"
    "E12" = "-2.51%"
    "D13" = "0.03453"
    "E13" = "-4.43%"
    "D14" = "0.09834"
    "E14" = "-0.13%"
    "D15" = "0.001422"
    "E15" = "2.37%"
    "D16" = "0.005747"
    "E16" = "0.14%"
    "D17" = "3.621"
    "E17" = "3.52%"
    "D18" = "2.969"
    "E18" = "2.21%"
    "D19" = "0.3439"
    "E19" = "-0.23%"
    "D20" = "0.1340"
    "E20" = "4.56%"
    "D21" = "4.896"
    "E21" = "-2.67%"
    "D22" = "0.2444"
    "E22" = "-0.81%"
    "D23" = "0.04455"
    "E23" = "-1.18%"
    "D24" = "0.001221"
    "E24" = "0.43%"
    "D25" = "0.004835"
    "E25" = "-0.40%"
    "D26" = "0.0001243"
    "E26" = "-0.48%"
    "D39" = "0.02025"
    "E39" = "4.78%"
    "D40" = "0.04935"
    "E40" = "4.01%"
    "D41" = "0.007709"
    "E41" = "1.97%"
    "D42" = "0.01006"
    "E42" = "5.32%"
    "D43" = "0.1378"
    "E43" = "3.59%"
    "D44" = "0.002106"
    "E44" = "-0.15%"
    "D45" = "0.01161"
    "E45" = "4.22%"
    "D46" = "0.00006449"
    "E46" = "2.38%"
    "E47" = "0.17%"
    "E49" = "-8.59%"
    "D50" = "0.00002103"
    "E50" = "0.17%"
    "D51" = "0.0002003"
    "E51" = "0.17%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
